# Realestate Update resale numbers 2024-01-13 21:41
# Append a new data row (row 55) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (A-D) must stay as literal text (matching the existing
# inlineStr cells), not get auto-converted into dates/times/numbers by
# Excel's "smart" value parsing. Temporarily force a text number format,
# assign the values, then clear formats again so no style index is left
# on the cells (matching the rest of the unstyled data rows).
$textRange = $ws.Range("A55:D55")
$textRange.NumberFormat = "@"

$ws.Range("A55").Value = "2024-01-13"
$ws.Range("B55").Value = "21:41:09"
$ws.Range("C55").Value = "Saturday"
$ws.Range("D55").Value = "01"

$textRange.ClearFormats()

# Numeric columns (E-T)
$ws.Range("E55").Value = 138585
$ws.Range("F55").Value = 142826
$ws.Range("G55").Value = 172084
$ws.Range("H55").Value = 148286
$ws.Range("I55").Value = -1
$ws.Range("J55").Value = 119761
$ws.Range("K55").Value = 225053
$ws.Range("L55").Value = 253549
$ws.Range("M55").Value = 185062
$ws.Range("N55").Value = 110445
$ws.Range("O55").Value = 41046
$ws.Range("P55").Value = 30915
$ws.Range("Q55").Value = 73148
$ws.Range("R55").Value = -1
$ws.Range("S55").Value = 42847
$ws.Range("T55").Value = -1
